$d = $word.ActiveDocument

# --- Paragraph 18: PM> add-migration ... command line ---
$xml18 = '<pkg:xmlData xmlns:pkg=''http://schemas.microsoft.com/office/2006/xmlPackage'' xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:p><w:r><w:rPr><w:rFonts w:ascii="Cascadia Mono" w:hAnsi="Cascadia Mono" w:cs="Cascadia Mono"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t xml:space="preserve">PM&gt; </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cascadia Mono" w:hAnsi="Cascadia Mono" w:cs="Cascadia Mono"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t xml:space="preserve">add-migration </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Cascadia Mono" w:hAnsi="Cascadia Mono" w:cs="Cascadia Mono"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>ShoppingCartMigrations.empty</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Cascadia Mono" w:hAnsi="Cascadia Mono" w:cs="Cascadia Mono"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>-test -verbose -</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Cascadia Mono" w:hAnsi="Cascadia Mono" w:cs="Cascadia Mono"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>startupproject</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Cascadia Mono" w:hAnsi="Cascadia Mono" w:cs="Cascadia Mono"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Cascadia Mono" w:hAnsi="Cascadia Mono" w:cs="Cascadia Mono"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>shoppingcartef</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Cascadia Mono" w:hAnsi="Cascadia Mono" w:cs="Cascadia Mono"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t xml:space="preserve"> -project </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Cascadia Mono" w:hAnsi="Cascadia Mono" w:cs="Cascadia Mono"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>shoppingcartmigrations</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></pkg:xmlData>'
$d.Paragraphs(18).Range.InsertXML($xml18)

# --- Paragraph 16: Add DesignTimeDbContextFactory to ShoppingCartEF ---
$xml16 = '<pkg:xmlData xmlns:pkg=''http://schemas.microsoft.com/office/2006/xmlPackage'' xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:p><w:r><w:t xml:space="preserve">Add </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>DesignTimeDbContextFactory</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ShoppingCartEF</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></pkg:xmlData>'
$d.Paragraphs(16).Range.InsertXML($xml16)

# --- Paragraph 14: OnConfigure ---
$xml14 = '<pkg:xmlData xmlns:pkg=''http://schemas.microsoft.com/office/2006/xmlPackage'' xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:p><w:r><w:tab/></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>OnConfigure</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> – revise “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>UserSqlServer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>” to include parameter for path to migrations project</w:t></w:r></w:p></pkg:xmlData>'
$d.Paragraphs(14).Range.InsertXML($xml14)

# --- Paragraph 13: Constructor ---
$xml13 = '<pkg:xmlData xmlns:pkg=''http://schemas.microsoft.com/office/2006/xmlPackage'' xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:p><w:r><w:tab/><w:t>Constructor</w:t></w:r><w:r><w:t xml:space="preserve"> – supply </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>DbContextOptions</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (pass to base) and connection string</w:t></w:r></w:p></pkg:xmlData>'
$d.Paragraphs(13).Range.InsertXML($xml13)

# --- Paragraph 12: Revise ShoppingCartDS ---
$xml12 = '<pkg:xmlData xmlns:pkg=''http://schemas.microsoft.com/office/2006/xmlPackage'' xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:p><w:r><w:t xml:space="preserve">Revise </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ShoppingCartDS</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p></pkg:xmlData>'
$d.Paragraphs(12).Range.InsertXML($xml12)

# --- Paragraph 11: Add AppSettings to ShoppingCartEF ---
$xml11 = '<pkg:xmlData xmlns:pkg=''http://schemas.microsoft.com/office/2006/xmlPackage'' xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:p><w:r><w:t xml:space="preserve">Add </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>AppSettings</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ShoppingCartEF</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></pkg:xmlData>'
$d.Paragraphs(11).Range.InsertXML($xml11)

# --- Paragraph 10: Add NuGet packages Microsoft.Extensions.Configuration.Json ---
$xml10 = '<pkg:xmlData xmlns:pkg=''http://schemas.microsoft.com/office/2006/xmlPackage'' xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:p><w:r><w:t xml:space="preserve">Add NuGet packages </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>Microsoft.Extensions.Configuration</w:t></w:r><w:r><w:t>.Json</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/></w:p></pkg:xmlData>'
$d.Paragraphs(10).Range.InsertXML($xml10)

# --- Paragraph 9: Add NuGet packages Microsoft.Extensions.Configuration ---
$xml9 = '<pkg:xmlData xmlns:pkg=''http://schemas.microsoft.com/office/2006/xmlPackage'' xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:p><w:r><w:t xml:space="preserve">Add NuGet packages </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>Microsoft.Extensions.Configuration</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/></w:p></pkg:xmlData>'
$d.Paragraphs(9).Range.InsertXML($xml9)

# --- Paragraph 8: Rename namespace -> expands into 4 paragraphs ---
$xml8 = '<pkg:xmlData xmlns:pkg=''http://schemas.microsoft.com/office/2006/xmlPackage'' xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:p><w:r><w:t>Rename namespace in migrations folder to match project name (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ShoppingCartMigrations</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Set Project Output Path of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ShoppingCartMigrations</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to the same folder as the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ShoppingCartEF</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> project </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Add reference to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ShoppingCartEF</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> from </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ShoppingCartMigrations</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p/></pkg:xmlData>'
$d.Paragraphs(8).Range.InsertXML($xml8)

# --- Paragraph 7 (old): "Add reference to ShoppingCartEF from ShoppingCartMigrations" -- remove (recreated later, after the new Set Project Output Path paragraph) ---
$d.Paragraphs(7).Range.Delete()

# --- Paragraph 5: Add new class library project (ShoppingCartMigrations) ---
$xml5 = '<pkg:xmlData xmlns:pkg=''http://schemas.microsoft.com/office/2006/xmlPackage'' xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:p><w:r><w:t>Add new class library project (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ShoppingCartMigrations</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p></pkg:xmlData>'
$d.Paragraphs(5).Range.InsertXML($xml5)

# --- Paragraph 3: Note: Must have at least one migration... ---
$xml3 = '<pkg:xmlData xmlns:pkg=''http://schemas.microsoft.com/office/2006/xmlPackage'' xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:p><w:r><w:t xml:space="preserve">Note: Must have at least one migration in project with data context before moving – this ensures the “Add-Migration” command can find the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dbContext</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p></pkg:xmlData>'
$d.Paragraphs(3).Range.InsertXML($xml3)

